$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: the product previously on row 18 ("معجون اسنان بارودونتكس
#     50مل") is replaced by a different product, "فرش  فوكس عرض", with an
#     updated quantity text and price.
$ws.Range("B18").Value = "فرش  فوكس عرض"
$ws.Range("H18").Value = "0:0"
$ws.Range("L18").Value = 40

# --- Step 2: insert a new row after row 18 for the product that used to
#     be there ("معجون اسنان بارودونتكس 50مل"), which now becomes row 19.
$ws.Rows.Item(19).Insert()
$ws.Range("A18:N18").Copy()
$ws.Range("A19:N19").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(19).RowHeight = 24.75
$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()

$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "معجون اسنان بارودونتكس 50مل"
$ws.Range("H19").Value = "6:0"
$ws.Range("L19").Value = 80
$ws.Range("N19").Value = "1:0"

# --- Step 3: insert another new row after row 19 for the brand-new
#     product, "معجون سيجنال 25 مل", which becomes row 20.
$ws.Rows.Item(20).Insert()
$ws.Range("A18:N18").Copy()
$ws.Range("A20:N20").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(20).RowHeight = 25.5
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()

$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "معجون سيجنال 25 مل"
$ws.Range("H20").Value = "1:0"
$ws.Range("L20").Value = 20
$ws.Range("N20").Value = "1:0"

$excel.CutCopyMode = 0

# --- Step 4: the totals row (previously row 19, now pushed down to row
#     21) is updated to reflect the new/changed rows.
$ws.Range("K21").Value = 1333
